$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values. Per commit message: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals" -- the K values
# for several rows need to be updated to their newly calculated values.

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 7
$ws.Range("G8").Value = 2
$ws.Range("G10").Value = 1
